# Regenerate merged AHB files
#
# - Rename the comparison-column headers from the "_old"/"_new" naming
#   scheme to the version-specific "_FV2304"/"_FV2310" naming scheme
#   (columns A-J = FV2304 side, K = diff, L-U = FV2310 side).
# - Turn the data range A1:U70 into a native Excel Table ("Table1") with
#   an AutoFilter on the header row.
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2304", "Segmentgruppe_FV2304", "Segment_FV2304", "Datenelement_FV2304", "Segment ID_FV2304",
    "Code_FV2304", "Qualifier_FV2304", "Beschreibung_FV2304", "Bedingungsausdruck_FV2304", "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310", "Segmentgruppe_FV2310", "Segment_FV2310", "Datenelement_FV2310", "Segment ID_FV2310",
    "Code_FV2310", "Qualifier_FV2310", "Beschreibung_FV2310", "Bedingungsausdruck_FV2310", "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert A1:U70 into a native Excel table with an autofilter on the header row
$tableRange = $ws.Range("A1:U70")
$lo = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Table1"

# Freeze the top (header) row
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
